$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix home/away ordering: swap columns F:V between paired rows ---
# Row 8
$ws.Range("F8").Value = 'AZ Alkmaar'
$ws.Range("G8").Value = 5.0
$ws.Range("H8").Value = 'G.A. Eagles'
$ws.Range("I8").Value = 1.0
$ws.Range("J8").Value = 1.37
$ws.Range("K8").Value = '07/07/2023 11:12'
$ws.Range("L8").Value = 1.54
$ws.Range("M8").Value = '13/08/2023 14:27'
$ws.Range("N8").Value = 5.51
$ws.Range("O8").Value = '07/07/2023 11:12'
$ws.Range("P8").Value = 4.41
$ws.Range("Q8").Value = '13/08/2023 14:26'
$ws.Range("R8").Value = 7.89
$ws.Range("S8").Value = '07/07/2023 11:12'
$ws.Range("T8").Value = 6.35
$ws.Range("U8").Value = '13/08/2023 14:26'
$ws.Range("V8").Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/az-alkmaar-g-a-eagles/SrIpNFQ5/'

# Row 9
$ws.Range("F9").Value = 'Feyenoord'
$ws.Range("G9").Value = 0.0
$ws.Range("H9").Value = 'Sittard'
$ws.Range("I9").Value = 0.0
$ws.Range("J9").Value = 1.22
$ws.Range("K9").Value = '07/07/2023 11:12'
$ws.Range("L9").Value = 1.19
$ws.Range("M9").Value = '13/08/2023 14:25'
$ws.Range("N9").Value = 7.61
$ws.Range("O9").Value = '07/07/2023 11:12'
$ws.Range("P9").Value = 7.38
$ws.Range("Q9").Value = '13/08/2023 14:29'
$ws.Range("R9").Value = 12.86
$ws.Range("S9").Value = '07/07/2023 11:12'
$ws.Range("T9").Value = 16.76
$ws.Range("U9").Value = '13/08/2023 14:29'
$ws.Range("V9").Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/feyenoord-sittard/6ZGlMZuC/'

# Row 17
$ws.Range("F17").Value = 'Sparta Rotterdam'
$ws.Range("G17").Value = 2.0
$ws.Range("H17").Value = 'Feyenoord'
$ws.Range("I17").Value = 2.0
$ws.Range("J17").Value = 5.0
$ws.Range("K17").Value = '13/08/2023 14:42'
$ws.Range("L17").Value = 4.94
$ws.Range("M17").Value = '20/08/2023 14:28'
$ws.Range("N17").Value = 4.52
$ws.Range("O17").Value = '13/08/2023 14:42'
$ws.Range("P17").Value = 4.35
$ws.Range("Q17").Value = '20/08/2023 14:29'
$ws.Range("R17").Value = 1.61
$ws.Range("S17").Value = '13/08/2023 14:42'
$ws.Range("T17").Value = 1.67
$ws.Range("U17").Value = '20/08/2023 14:28'
$ws.Range("V17").Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/sparta-rotterdam-feyenoord/21WSZhAB/'

# Row 18
$ws.Range("F18").Value = 'Twente'
$ws.Range("G18").Value = 3.0
$ws.Range("H18").Value = 'Zwolle'
$ws.Range("I18").Value = 1.0
$ws.Range("J18").Value = 1.35
$ws.Range("K18").Value = '13/08/2023 17:12'
$ws.Range("L18").Value = 1.44
$ws.Range("M18").Value = '20/08/2023 14:21'
$ws.Range("N18").Value = 5.56
$ws.Range("O18").Value = '13/08/2023 17:12'
$ws.Range("P18").Value = 5.21
$ws.Range("Q18").Value = '20/08/2023 14:28'
$ws.Range("R18").Value = 8.79
$ws.Range("S18").Value = '13/08/2023 17:12'
$ws.Range("T18").Value = 6.95
$ws.Range("U18").Value = '20/08/2023 14:29'
$ws.Range("V18").Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/twente-zwolle/KxsOzZf5/'

# Row 40
$ws.Range("F40").Value = 'Excelsior'
$ws.Range("G40").Value = 0.0
$ws.Range("H40").Value = 'Almere City'
$ws.Range("I40").Value = 0.0
$ws.Range("J40").Value = 1.93
$ws.Range("K40").Value = '04/09/2023 08:43'
$ws.Range("L40").Value = 2.04
$ws.Range("M40").Value = '17/09/2023 14:29'
$ws.Range("N40").Value = 3.96
$ws.Range("O40").Value = '04/09/2023 08:43'
$ws.Range("P40").Value = 3.95
$ws.Range("Q40").Value = '17/09/2023 14:28'
$ws.Range("R40").Value = 3.8
$ws.Range("S40").Value = '04/09/2023 08:43'
$ws.Range("T40").Value = 3.49
$ws.Range("U40").Value = '17/09/2023 14:29'
$ws.Range("V40").Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/excelsior-almere-city/fiHLuLx2/'

# Row 41
$ws.Range("F41").Value = 'Twente'
$ws.Range("G41").Value = 3.0
$ws.Range("H41").Value = 'Ajax'
$ws.Range("I41").Value = 1.0
$ws.Range("J41").Value = 3.43
$ws.Range("K41").Value = '04/09/2023 08:43'
$ws.Range("L41").Value = 2.47
$ws.Range("M41").Value = '17/09/2023 14:29'
$ws.Range("N41").Value = 4.07
$ws.Range("O41").Value = '04/09/2023 08:43'
$ws.Range("P41").Value = 3.6
$ws.Range("Q41").Value = '17/09/2023 14:28'
$ws.Range("R41").Value = 2.02
$ws.Range("S41").Value = '04/09/2023 08:43'
$ws.Range("T41").Value = 2.9
$ws.Range("U41").Value = '17/09/2023 14:29'
$ws.Range("V41").Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/twente-ajax/dde444F2/'

# Row 49
$ws.Range("F49").Value = 'Zwolle'
$ws.Range("G49").Value = 0.0
$ws.Range("H49").Value = 'AZ Alkmaar'
$ws.Range("I49").Value = 3.0
$ws.Range("J49").Value = 4.72
$ws.Range("K49").Value = '17/09/2023 16:13'
$ws.Range("L49").Value = 4.63
$ws.Range("M49").Value = '24/09/2023 16:39'
$ws.Range("N49").Value = 4.32
$ws.Range("O49").Value = '17/09/2023 16:13'
$ws.Range("P49").Value = 4.03
$ws.Range("Q49").Value = '24/09/2023 16:39'
$ws.Range("R49").Value = 1.67
$ws.Range("S49").Value = '17/09/2023 16:13'
$ws.Range("T49").Value = 1.77
$ws.Range("U49").Value = '24/09/2023 16:39'
$ws.Range("V49").Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/zwolle-az-alkmaar/ATQ3hbM7/'

# Row 50
$ws.Range("F50").Value = 'Waalwijk'
$ws.Range("G50").Value = 1.0
$ws.Range("H50").Value = 'Twente'
$ws.Range("I50").Value = 0.0
$ws.Range("J50").Value = 4.49
$ws.Range("K50").Value = '17/09/2023 13:43'
$ws.Range("L50").Value = 6.03
$ws.Range("M50").Value = '24/09/2023 16:43'
$ws.Range("N50").Value = 4.43
$ws.Range("O50").Value = '17/09/2023 13:43'
$ws.Range("P50").Value = 4.72
$ws.Range("Q50").Value = '24/09/2023 16:44'
$ws.Range("R50").Value = 1.68
$ws.Range("S50").Value = '17/09/2023 13:43'
$ws.Range("T50").Value = 1.53
$ws.Range("U50").Value = '24/09/2023 16:38'
$ws.Range("V50").Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/waalwijk-twente/OWEoHsa8/'

# Row 56
$ws.Range("F56").Value = 'Utrecht'
$ws.Range("G56").Value = 0.0
$ws.Range("H56").Value = 'Almere City'
$ws.Range("I56").Value = 2.0
$ws.Range("J56").Value = 1.51
$ws.Range("K56").Value = '23/09/2023 19:12'
$ws.Range("L56").Value = 1.57
$ws.Range("M56").Value = '30/09/2023 18:44'
$ws.Range("N56").Value = 4.75
$ws.Range("O56").Value = '23/09/2023 19:12'
$ws.Range("P56").Value = 4.4
$ws.Range("Q56").Value = '30/09/2023 18:44'
$ws.Range("R56").Value = 5.79
$ws.Range("S56").Value = '23/09/2023 19:12'
$ws.Range("T56").Value = 5.93
$ws.Range("U56").Value = '30/09/2023 18:44'
$ws.Range("V56").Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/utrecht-almere-city/dv2Y7wMf/'

# Row 57
$ws.Range("F57").Value = 'PSV'
$ws.Range("G57").Value = 3.0
$ws.Range("H57").Value = 'FC Volendam'
$ws.Range("I57").Value = 1.0
$ws.Range("J57").Value = 1.06
$ws.Range("K57").Value = '27/09/2023 18:12'
$ws.Range("L57").Value = 1.04
$ws.Range("M57").Value = '30/09/2023 14:21'
$ws.Range("N57").Value = 17.77
$ws.Range("O57").Value = '27/09/2023 18:12'
$ws.Range("P57").Value = 25.88
$ws.Range("Q57").Value = '30/09/2023 18:44'
$ws.Range("R57").Value = 24.06
$ws.Range("S57").Value = '27/09/2023 18:12'
$ws.Range("T57").Value = 42.53
$ws.Range("U57").Value = '30/09/2023 18:44'
$ws.Range("V57").Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/psv-fc-volendam/EFCT8J6l/'

# Row 83
$ws.Range("F83").Value = 'Sparta Rotterdam'
$ws.Range("G83").Value = 2.0
$ws.Range("H83").Value = 'Waalwijk'
$ws.Range("I83").Value = 0.0
$ws.Range("J83").Value = 1.75
$ws.Range("K83").Value = '22/10/2023 17:12'
$ws.Range("L83").Value = 1.78
$ws.Range("M83").Value = '28/10/2023 19:35'
$ws.Range("N83").Value = 3.96
$ws.Range("O83").Value = '22/10/2023 17:12'
$ws.Range("P83").Value = 3.95
$ws.Range("Q83").Value = '28/10/2023 19:37'
$ws.Range("R83").Value = 4.62
$ws.Range("S83").Value = '22/10/2023 17:12'
$ws.Range("T83").Value = 4.63
$ws.Range("U83").Value = '28/10/2023 19:35'
$ws.Range("V83").Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/sparta-rotterdam-waalwijk/YHIQY6QA/'

# Row 84
$ws.Range("F84").Value = 'Almere City'
$ws.Range("G84").Value = 0.0
$ws.Range("H84").Value = 'G.A. Eagles'
$ws.Range("I84").Value = 0.0
$ws.Range("J84").Value = 2.51
$ws.Range("K84").Value = '22/10/2023 17:12'
$ws.Range("L84").Value = 3.13
$ws.Range("M84").Value = '28/10/2023 19:58'
$ws.Range("N84").Value = 3.48
$ws.Range("O84").Value = '22/10/2023 17:12'
$ws.Range("P84").Value = 3.57
$ws.Range("Q84").Value = '28/10/2023 19:58'
$ws.Range("R84").Value = 2.89
$ws.Range("S84").Value = '22/10/2023 17:12'
$ws.Range("T84").Value = 2.34
$ws.Range("U84").Value = '28/10/2023 19:58'
$ws.Range("V84").Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/almere-city-g-a-eagles/xjPHzTeb/'

# --- Append two new match rows (97, 98), copying formatting from row 96 ---
$ws.Range("A96:V96").Copy()
$ws.Range("A97:V98").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Row 97
$ws.Range("A97").Value = 96
$ws.Range("B97").Value = 'netherlands'
$ws.Range("C97").Value = 'eredivisie'
$ws.Range("D97").Value = '2023-2024'
$ws.Range("E97").Value = 45235.69791666666
$ws.Range("F97").Value = 'Ajax'
$ws.Range("G97").Value = 4
$ws.Range("H97").Value = 'Heerenveen'
$ws.Range("I97").Value = 1
$ws.Range("J97").Value = 1.62
$ws.Range("K97").Value = '02/11/2023 20:12'
$ws.Range("L97").Value = 1.41
$ws.Range("M97").Value = '05/11/2023 16:44'
$ws.Range("N97").Value = 4.85
$ws.Range("O97").Value = '02/11/2023 20:12'
$ws.Range("P97").Value = 5.6
$ws.Range("Q97").Value = '05/11/2023 16:44'
$ws.Range("R97").Value = 4.55
$ws.Range("S97").Value = '02/11/2023 20:12'
$ws.Range("T97").Value = 7.02
$ws.Range("U97").Value = '05/11/2023 16:44'
$ws.Range("V97").Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/ajax-heerenveen/jgkikSAi/'

# Row 98
$ws.Range("A98").Value = 97
$ws.Range("B98").Value = 'netherlands'
$ws.Range("C98").Value = 'eredivisie'
$ws.Range("D98").Value = '2023-2024'
$ws.Range("E98").Value = 45235.69791666666
$ws.Range("F98").Value = 'Sparta Rotterdam'
$ws.Range("G98").Value = 1
$ws.Range("H98").Value = 'Almere City'
$ws.Range("I98").Value = 2
$ws.Range("J98").Value = 2.01
$ws.Range("K98").Value = '28/10/2023 20:13'
$ws.Range("L98").Value = 1.75
$ws.Range("M98").Value = '05/11/2023 16:36'
$ws.Range("N98").Value = 3.74
$ws.Range("O98").Value = '28/10/2023 20:13'
$ws.Range("P98").Value = 3.98
$ws.Range("Q98").Value = '05/11/2023 16:36'
$ws.Range("R98").Value = 3.73
$ws.Range("S98").Value = '28/10/2023 20:13'
$ws.Range("T98").Value = 4.79
$ws.Range("U98").Value = '05/11/2023 16:36'
$ws.Range("V98").Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/sparta-rotterdam-almere-city/z7ruv22j/'

